# Update database and change read_price algorithm
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Update the "twelve months ended" period headers (shift one year forward) ---
# Row 8 header (expense table)
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"

# Row 24 header (personnel table)
$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Shift the yearly data columns left by one, add the new 1401/12 column ---

# Row 10 - هزینه حمل و نقل و انتقال (unchanged zeros)
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0

# Row 11 - هزینه خدمات پس از فروش (unchanged zeros)
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0

# Row 12 - حق العمل و کمیسیون فروش (unchanged zeros)
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

# Row 13 - هزینه تبلیغات (unchanged zeros)
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

# Row 14 - هزینه مواد مصرفی (unchanged zeros)
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0

# Row 15 - هزینه انرژی (آب، برق، گاز و سوخت)
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 384
$ws.Range("I15").Value = 396

# Row 16 - هزینه استهلاک
$ws.Range("E16").Value = 549
$ws.Range("F16").Value = 788
$ws.Range("G16").Value = 1850
$ws.Range("H16").Value = 2039
$ws.Range("I16").Value = 1503

# Row 17 - هزینه حقوق و دستمزد
$ws.Range("E17").Value = 20179
$ws.Range("F17").Value = 32812
$ws.Range("G17").Value = 45274
$ws.Range("H17").Value = 66328
$ws.Range("I17").Value = 101089

# Row 18 - هزینه مطالبات مشکوک الوصول (unchanged zeros)
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

# Row 19 - سایر هزینه ها
$ws.Range("E19").Value = 32039
$ws.Range("F19").Value = 35554
$ws.Range("G19").Value = 49398
$ws.Range("H19").Value = 81137
$ws.Range("I19").Value = 128807

# Row 20 - جمع (total)
$ws.Range("E20").Value = 52767
$ws.Range("F20").Value = 69154
$ws.Range("G20").Value = 96522
$ws.Range("H20").Value = 149888
$ws.Range("I20").Value = 231795

# Row 26 - تعداد پرسنل غیر تولیدی شرکت
$ws.Range("E26").Value = 128
$ws.Range("F26").Value = 129
$ws.Range("G26").Value = 106
$ws.Range("H26").Value = 109
$ws.Range("I26").Value = 137

# Row 27 - تعداد پرسنل تولیدی شرکت
$ws.Range("E27").Value = 134
$ws.Range("F27").Value = 135
$ws.Range("G27").Value = 146
$ws.Range("H27").Value = 137
$ws.Range("I27").Value = 116
